$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 355, pushing existing rows 355..409 down to 356..410.
$ws.Rows("355").Insert()

# The row that used to be 355 is now row 356; duplicate it into the new
# row 355 as a starting point, then update the changed fields below.
$ws.Range("A356:R356").Copy()
$ws.Range("A355").PasteSpecial()

# Apply the new weekly entry's values.
$ws.Range("D355").Value2 = 44816
$ws.Range("J355").Value2 = 430
$ws.Range("K355").Value2 = 6000
$ws.Range("L355").Value2 = 6000
$ws.Range("M355").Value2 = 6000
$ws.Range("O355").Value2 = "Región de Arica y Parinacota"
$ws.Range("P355").Value2 = 6000
